$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the assigned person between rows 9-10 and rows 11-12
$ws.Range("D9").Value = "Đỗ Tiến Đạt"
$ws.Range("D10").Value = "Đỗ Tiến Đạt"
$ws.Range("D11").Value = "Lê Minh Hiếu"
$ws.Range("D12").Value = "Lê Minh Hiếu"

# Update the selected cell to D10
$ws.Range("D10").Select()
